$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38; this shifts the existing rows 38..152
# down to 39..153 (including formatting of the row above).
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly record.
$ws.Range("A38").Value = 4
$ws.Range("B38").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C38").Value = "Los Lagos"
$ws.Range("D38").Value = "2021-10-28"
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 100112017
$ws.Range("G38").Value = "Apio"
$ws.Range("H38").Value = "Americana (o)"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 35
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 11000
$ws.Range("M38").Value = 10429
$ws.Range("N38").Value = "`$/docena de matas"
$ws.Range("O38").Value = "Región de Coquimbo"
$ws.Range("P38").Value = 1738
$ws.Range("Q38").Value = 6
$ws.Range("R38").Value = "Hortaliza"
